# Auto-generated edit script applying scheduled market-price refresh
# to the Exodus_Profits workbook (per commit: "chore: update Sheets via scheduled runner")
$wb = $excel.ActiveWorkbook

# Sheet ALC, row 32 (Leve Item ID 5484)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 22833.334
$ws.Range("J32").Value = 17400
$ws.Range("L32").Value = 17400
$ws.Range("N32").Value = -18052

# Sheet ALC, row 51 (Leve Item ID 5486)
$ws.Range("H51").Value = 2790

# Sheet ALC, row 74 (Leve Item ID 5507)
$ws.Range("H74").Value = 4472.933
$ws.Range("I74").Value = 4314.923
$ws.Range("J74").Value = 5500
$ws.Range("K74").Value = 4314.923
$ws.Range("L74").Value = 5500
$ws.Range("M74").Value = -3378.923
$ws.Range("N74").Value = -7372

# Sheet ALC, row 77 (Leve Item ID 5507)
$ws.Range("H77").Value = 4472.933
$ws.Range("I77").Value = 4314.923
$ws.Range("J77").Value = 5500
$ws.Range("K77").Value = 21574.615
$ws.Range("L77").Value = 27500
$ws.Range("M77").Value = -16894.615
$ws.Range("N77").Value = -36860

# Sheet ALC, row 80 (Leve Item ID 12605)
$ws.Range("H80").Value = 1407.4736
$ws.Range("I80").Value = 295.125
$ws.Range("J80").Value = 2216.4546
$ws.Range("K80").Value = 885.375
$ws.Range("L80").Value = 6649.3638
$ws.Range("M80").Value = 112.625
$ws.Range("N80").Value = -8645.363799999999

# Sheet ALC, row 83 (Leve Item ID 12605)
$ws.Range("H83").Value = 1407.4736
$ws.Range("I83").Value = 295.125
$ws.Range("J83").Value = 2216.4546
$ws.Range("K83").Value = 2656.125
$ws.Range("L83").Value = 19948.0914
$ws.Range("M83").Value = 2335.875
$ws.Range("N83").Value = -29932.0914

# Sheet ALC, row 100 (Leve Item ID 19906)
$ws.Range("H100").Value = 1716.2727
$ws.Range("I100").Value = 1208.8889
$ws.Range("J100").Value = 3999.5
$ws.Range("K100").Value = 1208.8889
$ws.Range("L100").Value = 3999.5
$ws.Range("M100").Value = -667.8888999999999
$ws.Range("N100").Value = -5081.5

# Sheet ALC, row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 1587.6285
$ws.Range("I132").Value = 1619.6364
$ws.Range("J132").Value = 1059.5
$ws.Range("K132").Value = 4858.9092
$ws.Range("L132").Value = 3178.5
$ws.Range("M132").Value = -2328.9092
$ws.Range("N132").Value = -8238.5

# Sheet ALC, row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 520155.34
$ws.Range("I137").Value = 1421.4706
$ws.Range("K137").Value = 4264.4118
$ws.Range("M137").Value = -1714.4118

# Sheet ALC, row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 2983.2
$ws.Range("I138").Value = 2198.0833
$ws.Range("K138").Value = 6594.249899999999
$ws.Range("M138").Value = -1454.249899999999

# Sheet ARM, row 2 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1837.1578
$ws.Range("I2").Value = 1333.6428
$ws.Range("K2").Value = 1333.6428
$ws.Range("M2").Value = -1220.6428

# Sheet ARM, row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 36826.234
$ws.Range("I61").Value = 3532.875
$ws.Range("J61").Value = 169999.67
$ws.Range("K61").Value = 3532.875
$ws.Range("L61").Value = 169999.67
$ws.Range("M61").Value = -3320.875
$ws.Range("N61").Value = -170423.67

# Sheet ARM, row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 3215.3809
$ws.Range("I74").Value = 1822.5625
$ws.Range("J74").Value = 7672.4
$ws.Range("K74").Value = 1822.5625
$ws.Range("L74").Value = 7672.4
$ws.Range("M74").Value = -948.5625
$ws.Range("N74").Value = -9420.4

# Sheet ARM, row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 3215.3809
$ws.Range("I77").Value = 1822.5625
$ws.Range("J77").Value = 7672.4
$ws.Range("K77").Value = 9112.8125
$ws.Range("L77").Value = 38362
$ws.Range("M77").Value = -4744.8125
$ws.Range("N77").Value = -47098

# Sheet ARM, row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 1837.1578
$ws.Range("I116").Value = 1333.6428
$ws.Range("K116").Value = 1333.6428
$ws.Range("M116").Value = 960.3571999999999

# Sheet ARM, row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 2940.4583
$ws.Range("I132").Value = 2385.0698
$ws.Range("J132").Value = 7716.8
$ws.Range("K132").Value = 7155.209400000001
$ws.Range("L132").Value = 23150.4
$ws.Range("M132").Value = -4625.209400000001
$ws.Range("N132").Value = -28210.4

# Sheet ARM, row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 36826.234
$ws.Range("I136").Value = 3532.875
$ws.Range("J136").Value = 169999.67
$ws.Range("K136").Value = 10598.625
$ws.Range("L136").Value = 509999.01
$ws.Range("M136").Value = -8048.625
$ws.Range("N136").Value = -515099.01

# Sheet BSM, row 3 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1837.1578
$ws.Range("I3").Value = 1333.6428
$ws.Range("K3").Value = 1333.6428
$ws.Range("M3").Value = -1219.6428

# Sheet BSM, row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 6520
$ws.Range("I86").Value = 3025.4
$ws.Range("K86").Value = 3025.4
$ws.Range("M86").Value = -1902.4

# Sheet BSM, row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 6520
$ws.Range("I89").Value = 3025.4
$ws.Range("K89").Value = 15127
$ws.Range("M89").Value = -9511

# Sheet CRP, row 31 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2252.5833
$ws.Range("I31").Value = 1438.3077
$ws.Range("J31").Value = 3214.9092
$ws.Range("K31").Value = 1438.3077
$ws.Range("L31").Value = 3214.9092
$ws.Range("M31").Value = -1143.3077
$ws.Range("N31").Value = -3804.9092

# Sheet CRP, row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 2252.5833
$ws.Range("I34").Value = 1438.3077
$ws.Range("J34").Value = 3214.9092
$ws.Range("K34").Value = 1438.3077
$ws.Range("L34").Value = 3214.9092
$ws.Range("M34").Value = -1236.3077
$ws.Range("N34").Value = -3618.9092

# Sheet CRP, row 64 (Leve Item ID 10610)
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# Sheet CRP, row 67 (Leve Item ID 10610)
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# Sheet CRP, row 86 (Leve Item ID 12584)
$ws.Range("H86").Value = 2112284.8
$ws.Range("J86").Value = 15435.6
$ws.Range("L86").Value = 15435.6
$ws.Range("N86").Value = -17681.6

# Sheet CRP, row 89 (Leve Item ID 12584)
$ws.Range("H89").Value = 2112284.8
$ws.Range("J89").Value = 15435.6
$ws.Range("L89").Value = 77178
$ws.Range("N89").Value = -88410

# Sheet CRP, row 122 (Leve Item ID 36196)
$ws.Range("H122").Value = 1923.4286
$ws.Range("I122").Value = 1382.7333
$ws.Range("J122").Value = 2328.95
$ws.Range("K122").Value = 4148.199900000001
$ws.Range("L122").Value = 6986.849999999999
$ws.Range("M122").Value = -1698.199900000001
$ws.Range("N122").Value = -11886.85

# Sheet CRP, row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 1085220.6
$ws.Range("I132").Value = 1265651.9
$ws.Range("J132").Value = 2632.6667
$ws.Range("K132").Value = 3796955.7
$ws.Range("L132").Value = 7898.000100000001
$ws.Range("M132").Value = -3794425.7
$ws.Range("N132").Value = -12958.0001

# Sheet CRP, row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 3814481.5
$ws.Range("I134").Value = 7145759.5
$ws.Range("J134").Value = 113061.89
$ws.Range("K134").Value = 21437278.5
$ws.Range("L134").Value = 339185.67
$ws.Range("M134").Value = -21434743.5
$ws.Range("N134").Value = -344255.67

# Sheet CUL, row 12 (Leve Item ID 4854)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 45.63158
$ws.Range("I12").Value = 31
$ws.Range("J12").Value = 58.8
$ws.Range("K12").Value = 93
$ws.Range("L12").Value = 176.4
$ws.Range("M12").Value = 80
$ws.Range("N12").Value = -522.4

# Sheet CUL, row 59 (Leve Item ID 4694)
$ws.Range("H59").Value = 3221.4167
$ws.Range("I59").Value = 1900.3334
$ws.Range("J59").Value = 4542.5
$ws.Range("K59").Value = 5701.0002
$ws.Range("L59").Value = 13627.5
$ws.Range("M59").Value = -5161.0002
$ws.Range("N59").Value = -14707.5

# Sheet CUL, row 132 (Leve Item ID 43972)
$ws.Range("H132").Value = 798
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 798
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 7182
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -12242

# Sheet GSM, row 19 (Leve Item ID 2668)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 840
$ws.Range("I19").Value = 840
$ws.Range("K19").Value = 840
$ws.Range("M19").Value = -552

# Sheet GSM, row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 1859.3684
$ws.Range("I102").Value = 1812.9445
$ws.Range("J102").Value = 2695
$ws.Range("K102").Value = 1812.9445
$ws.Range("L102").Value = 2695
$ws.Range("M102").Value = -190.9445000000001
$ws.Range("N102").Value = -5939

# Sheet GSM, row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 7370.381
$ws.Range("I132").Value = 5633.625
$ws.Range("J132").Value = 12928
$ws.Range("K132").Value = 16900.875
$ws.Range("L132").Value = 38784
$ws.Range("M132").Value = -14370.875
$ws.Range("N132").Value = -43844

# Sheet LTW, row 40 (Leve Item ID 36248)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2178500
$ws.Range("I40").Value = 48866.047
$ws.Range("J40").Value = 13891488
$ws.Range("K40").Value = 48866.047
$ws.Range("L40").Value = 13891488
$ws.Range("M40").Value = -48730.047
$ws.Range("N40").Value = -13891760

# Sheet LTW, row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 1325.2
$ws.Range("I46").Value = 1242.3334
$ws.Range("J46").Value = 1449.5
$ws.Range("K46").Value = 1242.3334
$ws.Range("L46").Value = 1449.5
$ws.Range("M46").Value = -1054.3334
$ws.Range("N46").Value = -1825.5

# Sheet LTW, row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 1729
$ws.Range("I55").Value = 2071.2856
$ws.Range("J55").Value = 1249.8
$ws.Range("K55").Value = 2071.2856
$ws.Range("L55").Value = 1249.8
$ws.Range("M55").Value = -1898.2856
$ws.Range("N55").Value = -1595.8

# Sheet WVR, row 2 (Leve Item ID 3307)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 5492939.5
$ws.Range("J2").Value = 48002
$ws.Range("L2").Value = 48002
$ws.Range("N2").Value = -48226

# Sheet WVR, row 13 (Leve Item ID 3008)
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()

# Sheet WVR, row 16 (Leve Item ID 26304)
$ws.Range("H16").Value = 70420
$ws.Range("J16").Value = 70420
$ws.Range("L16").Value = 70420
$ws.Range("N16").Value = -71004

# Sheet WVR, row 81 (Leve Item ID 12596)
$ws.Range("H81").Value = 31271.312
$ws.Range("I81").Value = 3122.6365
$ws.Range("J81").Value = 93198.39999999999
$ws.Range("K81").Value = 6245.273
$ws.Range("L81").Value = 186396.8
$ws.Range("M81").Value = -5184.273
$ws.Range("N81").Value = -188518.8

# Sheet WVR, row 84 (Leve Item ID 12596)
$ws.Range("H84").Value = 31271.312
$ws.Range("I84").Value = 3122.6365
$ws.Range("J84").Value = 93198.39999999999
$ws.Range("K84").Value = 31226.365
$ws.Range("L84").Value = 931984
$ws.Range("M84").Value = -25922.365
$ws.Range("N84").Value = -942592

# Sheet WVR, row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 2224.422
$ws.Range("I132").Value = 2136.5
$ws.Range("J132").Value = 2576.111
$ws.Range("K132").Value = 6409.5
$ws.Range("L132").Value = 7728.333
$ws.Range("M132").Value = -3879.5
$ws.Range("N132").Value = -12788.333
